# Weekly fruit/vegetable price refresh: the per-row "measurement" columns
# (Fecha, Calidad, Volumen, Precio minimo/maximo/ponderado, Unidad de
# comercializacion, Origen, Precio $/Kg, Kg / unidad) get reshuffled across
# the data rows (2-25), while the descriptive columns (A,B,C,E,F,G,H,I,J,K)
# stay put. Build the new row->row mapping, snapshot all source values
# first, then write them all back so the rewrite is safe regardless of
# write order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (i.e. row $target ends up with the data that
# used to live in row $source)
$map = @{
    2  = 19
    3  = 12
    4  = 10
    5  = 20
    6  = 21
    7  = 11
    8  = 2
    9  = 14
    10 = 15
    11 = 4
    12 = 22
    13 = 23
    14 = 8
    15 = 9
    16 = 3
    17 = 16
    18 = 5
    19 = 24
    20 = 25
    21 = 18
    22 = 6
    23 = 7
    24 = 13
    25 = 17
}

# Columns that move together as a row's "measurement" payload.
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)  # D, L, M, N, O, P, Q, R, S, T

# Snapshot every source row's values before writing anything.
$snapshot = @{}
foreach ($srcRow in $map.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowVals = @{}
        foreach ($col in $cols) {
            $rowVals[$col] = $ws.Cells.Item($srcRow, $col).Value()
        }
        $snapshot[$srcRow] = $rowVals
    }
}

# Now write the snapshotted values into their destination rows.
foreach ($targetRow in $map.Keys) {
    $srcRow = $map[$targetRow]
    $rowVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($targetRow, $col).Value = $rowVals[$col]
    }
}
